$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85
$ws.Cells.Item($row, 1).Value = "GJ4F48"
$ws.Cells.Item($row, 2).Value = "Rodillo de entrega de papel de cubierta superior de fusor para HP"
$ws.Cells.Item($row, 3).Value = "M101 M102 M103 M104 M106 M129 M130 M131 M132 M133 M134 M203 M227 M230"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 50000
$ws.Cells.Item($row, 6).Value = 8
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E85-D85)*G85"
$ws.Cells.Item($row, 9).Formula = "=D85*F85"
$ws.Cells.Item($row, 10).Value = 0
